$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) so the new
# headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows 2-43: team record values (Wins=86, Losses=76, Ties=0)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
